$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-19 14:36:42"
$wsZh.Range("H4").Value = "2016-03-19 14:37:06"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-19 14:36:45"
$wsDe.Range("H4").Value = "2016-03-19 14:37:11"
